$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values per row (columns G,H,I,J,M,N,O,P,Q,R,S,T)
$data = @{
    2  = @{ G=2.727484333333333;  H=8.182453000000001;  I=0.03096049453772388; J=0.03096049453772388; M=0.110552;           N=0.331656;          O=0.01126249561724847; P=0.01126249561724847; Q=0.3015288480186666;  R=2.713759632168;     S=0.0003486924340389605; T=0.0003486924340389605 }
    3  = @{ G=2.727484333333333;  H=8.182453000000001;  I=0.03096049453772388; J=0.03096049453772388; M=9.012070666666666;  N=27.036212;          O=0.9181055646724333;  P=0.9181055646724334;  Q=24.58028155422622;   R=221.222533988036;    S=0.02842500232009477;   T=0.02842500232009478 }
    4  = @{ G=2.727484333333333;  H=8.182453000000001;  I=0.03096049453772388; J=0.03096049453772388; M=0.6933189999999999; N=2.079957;           O=0.07063193971031816; P=0.07063193971031817; Q=1.891016710502333;   R=17.019150394521;     S=0.002186799783590148;  T=0.002186799783590149 }
    5  = @{ I=0.5986009007423507; J=0.5986009007423507; M=0.110552;           N=0.331656;          O=0.01126249561724847; P=0.01126249561724847; Q=5.829862950149333;   R=52.468766551344;     S=0.006741740021091713;  T=0.006741740021091714 }
    6  = @{ I=0.5986009007423507; J=0.5986009007423507; O=0.9181055646724333;  P=0.9181055646724334;  Q=475.2436580408097;   R=4277.192922367288;  S=0.5495788179894832;    T=0.5495788179894832 }
    7  = @{ I=0.5986009007423507; J=0.5986009007423507; M=0.6933189999999999; N=2.079957;           O=0.07063193971031816; P=0.07063193971031817; Q=36.56157057976866;   R=329.054135217918;    S=0.04228034273177586;   T=0.04228034273177587 }
    8  = @{ G=32.63402300000001;  H=97.90206900000001;  I=0.3704386047199253;  J=0.3704386047199253;  M=0.110552;           N=0.331656;          O=0.01126249561724847; P=0.01126249561724847; Q=3.607756510696001;   R=32.469808596264;     S=0.004172063162117798;  T=0.004172063162117799 }
    9  = @{ G=32.63402300000001;  H=97.90206900000001;  I=0.3704386047199253;  J=0.3704386047199253;  O=0.9181055646724333;  P=0.9181055646724334;  Q=294.1001214136254;   R=2646.901092722628;   S=0.3401017443628553;    T=0.3401017443628554 }
    10 = @{ G=32.63402300000001;  H=97.90206900000001;  I=0.3704386047199253;  J=0.3704386047199253;  M=0.6933189999999999; N=2.079957;           O=0.07063193971031816; P=0.07063193971031817; Q=22.625788192337;     R=203.632093731033;    S=0.02616479719495214;   T=0.02616479719495215 }
}

foreach ($rowNum in $data.Keys) {
    $row = $data[$rowNum]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$rowNum").Value = $row[$col]
    }
}
